$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 691.8570999999999
$ws.Range("I31").Value = 691.8570999999999
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2075.5713
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1845.5713
$ws.Range("N31").ClearContents()

$ws.Range("H33").Value = 2854.95
$ws.Range("I33").Value = 2487.75
$ws.Range("J33").Value = 3405.75
$ws.Range("K33").Value = 2487.75
$ws.Range("L33").Value = 3405.75
$ws.Range("M33").Value = -2258.75
$ws.Range("N33").Value = -3863.75

$ws.Range("H64").Value = 4551.4736

$ws.Range("H67").Value = 4551.4736

$ws.Range("H98").Value = 55556890
$ws.Range("I98").Value = 62501136
$ws.Range("K98").Value = 62501136
$ws.Range("M98").Value = -62499638

$ws.Range("H111").Value = 3303.5
$ws.Range("I111").Value = 3571.3333
$ws.Range("J111").Value = 2500
$ws.Range("K111").Value = 10713.9999
$ws.Range("L111").Value = 7500
$ws.Range("M111").Value = -7646.999899999999
$ws.Range("N111").Value = -13634

$ws.Range("H122").Value = 55556890
$ws.Range("I122").Value = 62501136
$ws.Range("K122").Value = 187503408
$ws.Range("M122").Value = -187500958

$ws.Range("H136").Value = 47414.082
$ws.Range("J136").Value = 47414.082
$ws.Range("L136").Value = 47414.082
$ws.Range("N136").Value = -57614.082

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 876.75
$ws.Range("I2").Value = 752.25
$ws.Range("K2").Value = 752.25
$ws.Range("M2").Value = -639.25

$ws.Range("H32").Value = 17249384
$ws.Range("I32").Value = 27780064
$ws.Range("K32").Value = 27780064
$ws.Range("M32").Value = -27779777

$ws.Range("H61").Value = 62510416
$ws.Range("I61").Value = 62510000
$ws.Range("K61").Value = 62510000
$ws.Range("M61").Value = -62509788

$ws.Range("H74").Value = 21667756
$ws.Range("J74").Value = 5001005
$ws.Range("L74").Value = 5001005
$ws.Range("N74").Value = -5002753

$ws.Range("H77").Value = 21667756
$ws.Range("J77").Value = 5001005
$ws.Range("L77").Value = 25005025
$ws.Range("N77").Value = -25013761

$ws.Range("H88").Value = 3910.2222
$ws.Range("I88").Value = 2681.6667
$ws.Range("J88").Value = 4261.2383
$ws.Range("K88").Value = 2681.6667
$ws.Range("L88").Value = 4261.2383
$ws.Range("M88").Value = -2275.6667
$ws.Range("N88").Value = -5073.2383

$ws.Range("H91").Value = 3910.2222
$ws.Range("I91").Value = 2681.6667
$ws.Range("J91").Value = 4261.2383
$ws.Range("K91").Value = 2681.6667
$ws.Range("L91").Value = 4261.2383
$ws.Range("M91").Value = -1277.6667
$ws.Range("N91").Value = -7069.2383

$ws.Range("H116").Value = 876.75
$ws.Range("I116").Value = 752.25
$ws.Range("K116").Value = 752.25
$ws.Range("M116").Value = 1541.75

$ws.Range("H118").Value = 35000
$ws.Range("J118").Value = 35000
$ws.Range("L118").Value = 35000
$ws.Range("N118").Value = -38314

$ws.Range("H122").Value = 3494
$ws.Range("I122").Value = 1989.5
$ws.Range("J122").Value = 4998.5
$ws.Range("K122").Value = 5968.5
$ws.Range("L122").Value = 14995.5
$ws.Range("M122").Value = -3518.5
$ws.Range("N122").Value = -19895.5

$ws.Range("H132").Value = 2843.3125
$ws.Range("I132").Value = 2870.2903
$ws.Range("K132").Value = 8610.8709
$ws.Range("M132").Value = -6080.8709

$ws.Range("H136").Value = 62510416
$ws.Range("I136").Value = 62510000
$ws.Range("K136").Value = 187530000
$ws.Range("M136").Value = -187527450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 876.75
$ws.Range("I3").Value = 752.25
$ws.Range("K3").Value = 752.25
$ws.Range("M3").Value = -638.25

$ws.Range("H22").Value = 1099.25
$ws.Range("I22").Value = 822.5
$ws.Range("J22").Value = 1376
$ws.Range("K22").Value = 822.5
$ws.Range("L22").Value = 1376
$ws.Range("M22").Value = -649.5
$ws.Range("N22").Value = -1722

$ws.Range("H81").Value = 27273.166
$ws.Range("J81").Value = 27273.166
$ws.Range("L81").Value = 27273.166
$ws.Range("N81").Value = -29395.166

$ws.Range("H84").Value = 27273.166
$ws.Range("J84").Value = 27273.166
$ws.Range("L84").Value = 81819.49800000001
$ws.Range("N84").Value = -92427.49800000001

$ws.Range("H113").Value = 5465.6665
$ws.Range("I113").Value = 5465.6665
$ws.Range("K113").Value = 5465.6665
$ws.Range("M113").Value = -3295.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2578.1052
$ws.Range("I58").Value = 2227.5557
$ws.Range("J58").Value = 8888
$ws.Range("K58").Value = 2227.5557
$ws.Range("L58").Value = 8888
$ws.Range("M58").Value = -2024.5557
$ws.Range("N58").Value = -9294

$ws.Range("H76").Value = 6753.1665
$ws.Range("I76").Value = 6753.1665
$ws.Range("K76").Value = 6753.1665
$ws.Range("M76").Value = -6438.1665

$ws.Range("H79").Value = 6753.1665
$ws.Range("I79").Value = 6753.1665
$ws.Range("K79").Value = 6753.1665
$ws.Range("M79").Value = -5661.1665

$ws.Range("H99").Value = 3704.6
$ws.Range("I99").Value = 3528.7273
$ws.Range("K99").Value = 3528.7273
$ws.Range("M99").Value = -2030.7273

$ws.Range("H107").Value = 2008.96
$ws.Range("I107").Value = 866.7692
$ws.Range("K107").Value = 866.7692
$ws.Range("M107").Value = 1053.2308

$ws.Range("H122").Value = 1874.6428
$ws.Range("I122").Value = 1577.7273
$ws.Range("K122").Value = 4733.1819
$ws.Range("M122").Value = -2283.1819

$ws.Range("H126").Value = 3704.6
$ws.Range("I126").Value = 3528.7273
$ws.Range("K126").Value = 10586.1819
$ws.Range("M126").Value = -8116.1819

$ws.Range("H134").Value = 3060.8572
$ws.Range("I134").Value = 3060.8572
$ws.Range("K134").Value = 9182.571599999999
$ws.Range("M134").Value = -6647.571599999999

$ws.Range("H136").Value = 2578.1052
$ws.Range("I136").Value = 2227.5557
$ws.Range("J136").Value = 8888
$ws.Range("K136").Value = 6682.6671
$ws.Range("L136").Value = 26664
$ws.Range("M136").Value = -4132.6671
$ws.Range("N136").Value = -31764

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2255.2307
$ws.Range("I5").Value = 2260.6667
$ws.Range("K5").Value = 6782.000100000001
$ws.Range("M5").Value = -6670.000100000001

$ws.Range("H8").Value = 166.7
$ws.Range("I8").Value = 166.7
$ws.Range("K8").Value = 500.1
$ws.Range("M8").Value = -361.1

$ws.Range("H12").Value = 2437.5454
$ws.Range("I12").Value = 3327.5715
$ws.Range("K12").Value = 9982.7145
$ws.Range("M12").Value = -9809.7145

$ws.Range("H135").Value = 2255.2307
$ws.Range("I135").Value = 2260.6667
$ws.Range("K135").Value = 20346.0003
$ws.Range("M135").Value = -17811.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3122.1052
$ws.Range("I80").Value = 2458.1875
$ws.Range("K80").Value = 2458.1875
$ws.Range("M80").Value = -1460.1875

$ws.Range("H83").Value = 3122.1052
$ws.Range("I83").Value = 2458.1875
$ws.Range("K83").Value = 12290.9375
$ws.Range("M83").Value = -7298.9375

$ws.Range("H122").Value = 2073.4546
$ws.Range("I122").Value = 1867.7778
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 5603.3334
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -3153.3334
$ws.Range("N122").Value = -13897

$ws.Range("H132").Value = 34490268
$ws.Range("I132").Value = 47622428
$ws.Range("J132").Value = 18347.25
$ws.Range("K132").Value = 142867284
$ws.Range("L132").Value = 55041.75
$ws.Range("M132").Value = -142864754
$ws.Range("N132").Value = -60101.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 27000
$ws.Range("I11").Value = 27000
$ws.Range("K11").Value = 27000
$ws.Range("M11").Value = -26860

$ws.Range("H123").Value = 49330
$ws.Range("J123").Value = 49330
$ws.Range("L123").Value = 49330
$ws.Range("N123").Value = -59130

$ws.Range("H132").Value = 1266427.4
$ws.Range("I132").Value = 28201
$ws.Range("K132").Value = 84603
$ws.Range("M132").Value = -82073

$ws.Range("H133").Value = 68000
$ws.Range("J133").Value = 68000
$ws.Range("L133").Value = 68000
$ws.Range("N133").Value = -73060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2115.8948
$ws.Range("I122").Value = 2149.3438
$ws.Range("J122").Value = 1937.5
$ws.Range("K122").Value = 6448.0314
$ws.Range("L122").Value = 5812.5
$ws.Range("M122").Value = -3998.0314
$ws.Range("N122").Value = -10712.5

$ws.Range("H132").Value = 187502.61
$ws.Range("I132").Value = 2318.239
$ws.Range("K132").Value = 6954.717000000001
$ws.Range("M132").Value = -4424.717000000001

$ws.Range("H136").Value = 5275.25
$ws.Range("I136").Value = 6118.273
$ws.Range("K136").Value = 18354.819
$ws.Range("M136").Value = -15804.819
